$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.021.38"
$ws.Range("E2").Value = "  +0.55%  "
$ws.Range("D3").Value = "1.641.30"
$ws.Range("E3").Value = "  +0.23%  "
$ws.Range("E4").Value = "  -0.23%  "
$ws.Range("D5").Value = "'214.50"
$ws.Range("E5").Value = "  +0.13%  "
$ws.Range("D6").Value = "'0.5093"
$ws.Range("E6").Value = "  +1.52%  "
$ws.Range("E7").Value = "  -0.14%  "
$ws.Range("D8").Value = "'0.2563"
$ws.Range("E8").Value = "  +0.03%  "
$ws.Range("D9").Value = "'0.06356"
$ws.Range("E9").Value = "  -0.34%  "
$ws.Range("D10").Value = "'19.56"
$ws.Range("E10").Value = "  +0.32%  "
$ws.Range("D11").Value = "'0.07759"
$ws.Range("E11").Value = "  -0.34%  "
$ws.Range("E12").Value = "  +0.46%  "
$ws.Range("D13").Value = "1.643.94"
$ws.Range("E13").Value = "  +0.30%  "
$ws.Range("D14").Value = "'0.5439"
$ws.Range("E14").Value = "  +0.47%  "
$ws.Range("B15").Value = "ShibaInu"
$ws.Range("C15").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D15").Value = "0.0₅7718"
$ws.Range("E15").Value = "  -1.63%  "
$ws.Range("B16").Value = "Litecoin"
$ws.Range("C16").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D16").Value = "'64.18"
$ws.Range("E16").Value = "  -0.52%  "
$ws.Range("D17").Value = "26.020.51"
$ws.Range("E17").Value = "  +0.45%  "
$ws.Range("E18").Value = "  -0.18%  "
$ws.Range("D19").Value = "'197.19"
$ws.Range("E19").Value = "  -0.14%  "
$ws.Range("D20").Value = "'4.419"
$ws.Range("E20").Value = "  +1.00%  "
$ws.Range("D21").Value = "'9.922"
$ws.Range("E21").Value = "  +0.13%  "
$ws.Range("D22").Value = "'6.027"
$ws.Range("E22").Value = "  +1.16%  "
$ws.Range("E23").Value = "  -0.08%  "
$ws.Range("D24").Value = "'1.863"
$ws.Range("E24").Value = "  -0.47%  "
$ws.Range("D25").Value = "'140.65"
$ws.Range("E25").Value = "  +0.67%  "
$ws.Range("D26").Value = "'0.1195"
$ws.Range("E26").Value = "  +5.03%  "
$ws.Range("D27").Value = "'6.809"
$ws.Range("E27").Value = "  -0.13%  "
$ws.Range("D28").Value = "'15.53"
$ws.Range("E28").Value = "  -0.52%  "
$ws.Range("E29").Value = "  -0.38%  "
$ws.Range("D30").Value = "'0.04846"
$ws.Range("E30").Value = "  -0.55%  "
$ws.Range("D31").Value = "'3.250"
$ws.Range("E31").Value = "  +0.04%  "
$ws.Range("D32").Value = "'3.165"
$ws.Range("E32").Value = "  -0.60%  "
$ws.Range("D33").Value = "'1.525"
$ws.Range("E33").Value = "  -0.06%  "
$ws.Range("E34").Value = "  -0.15%  "
$ws.Range("D35").Value = "'0.8982"
$ws.Range("E35").Value = "  +1.28%  "
$ws.Range("D36").Value = "'2.576"
$ws.Range("E36").Value = "  -0.79%  "
$ws.Range("D37").Value = "1.141.17"
$ws.Range("E37").Value = "  +0.84%  "
$ws.Range("E38").Value = "  -1.29%  "
$ws.Range("E39").Value = "  +0.23%  "
$ws.Range("E40").Value = "  -0.11%  "
$ws.Range("D41").Value = "'2.537"
$ws.Range("E41").Value = "  -0.72%  "
$ws.Range("E42").Value = "  +6.81%  "
$ws.Range("E43").Value = "  -0.73%  "
$ws.Range("D44").Value = "'99.24"
$ws.Range("E44").Value = "  +0.02%  "
$ws.Range("D45").Value = "'5.379"
$ws.Range("E45").Value = "  -5.15%  "
$ws.Range("D46").Value = "1.779.40"
$ws.Range("E46").Value = "  +0.27%  "
$ws.Range("D47").Value = "'0.4523"
$ws.Range("E47").Value = "  +0.10%  "
$ws.Range("B48").Value = "Aave"
$ws.Range("C48").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D48").Value = "'54.82"
$ws.Range("E48").Value = "  -0.47%  "
$ws.Range("B49").Value = "Frax"
$ws.Range("C49").Value = "https://coinranking.com/coin/KfWtaeV1W+frax-frax"
$ws.Range("D49").Value = "'0.9976"
$ws.Range("E49").Value = "  -0.79%  "
$ws.Range("D50").Value = "'0.05053"
$ws.Range("E50").Value = "  -0.57%  "
$ws.Range("D51").Value = "'1.001"
$ws.Range("E51").Value = "  -0.21%  "
